$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 686.46155
$ws.Range("J19").Value = 708.2941
$ws.Range("L19").Value = 708.2941
$ws.Range("N19").Value = -1058.2941
$ws.Range("H106").Value = 22224204
$ws.Range("I106").Value = 27779504
$ws.Range("K106").Value = 27779504
$ws.Range("M106").Value = -27778873
$ws.Range("H111").Value = 15642.286
$ws.Range("J111").Value = 1998.6666
$ws.Range("L111").Value = 5995.9998
$ws.Range("N111").Value = -12129.9998
$ws.Range("H137").Value = 1423.9166
$ws.Range("I137").Value = 1055.6666
$ws.Range("J137").Value = 1792.1666
$ws.Range("K137").Value = 3166.9998
$ws.Range("L137").Value = 5376.4998
$ws.Range("M137").Value = -616.9998000000001
$ws.Range("N137").Value = -10476.4998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H61").Value = 1850.1111
$ws.Range("I61").Value = 1059.75
$ws.Range("K61").Value = 1059.75
$ws.Range("M61").Value = -847.75
$ws.Range("H132").Value = 2449.4517
$ws.Range("I132").Value = 1962.125
$ws.Range("J132").Value = 4120.2856
$ws.Range("K132").Value = 5886.375
$ws.Range("L132").Value = 12360.8568
$ws.Range("M132").Value = -3356.375
$ws.Range("N132").Value = -17420.8568
$ws.Range("H136").Value = 1850.1111
$ws.Range("I136").Value = 1059.75
$ws.Range("K136").Value = 3179.25
$ws.Range("M136").Value = -629.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1319.5
$ws.Range("I36").Value = 1319.5
$ws.Range("K36").Value = 1319.5
$ws.Range("M36").Value = -785.5
$ws.Range("H64").Value = 472.63635
$ws.Range("J64").Value = 505.83334
$ws.Range("L64").Value = 505.83334
$ws.Range("N64").Value = -955.83334
$ws.Range("H67").Value = 472.63635
$ws.Range("J67").Value = 505.83334
$ws.Range("L67").Value = 505.83334
$ws.Range("N67").Value = -2065.83334
$ws.Range("H86").Value = 27501.166
$ws.Range("I86").Value = 13000
$ws.Range("K86").Value = 13000
$ws.Range("M86").Value = -11877
$ws.Range("H89").Value = 27501.166
$ws.Range("I89").Value = 13000
$ws.Range("K89").Value = 65000
$ws.Range("M89").Value = -59384
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1800.4546
$ws.Range("I31").Value = 1122.7
$ws.Range("J31").Value = 2365.25
$ws.Range("K31").Value = 1122.7
$ws.Range("L31").Value = 2365.25
$ws.Range("M31").Value = -827.7
$ws.Range("N31").Value = -2955.25
$ws.Range("H34").Value = 1800.4546
$ws.Range("I34").Value = 1122.7
$ws.Range("J34").Value = 2365.25
$ws.Range("K34").Value = 1122.7
$ws.Range("L34").Value = 2365.25
$ws.Range("M34").Value = -920.7
$ws.Range("N34").Value = -2769.25
$ws.Range("H58").Value = 1142.9767
$ws.Range("I58").Value = 574.931
$ws.Range("J58").Value = 2319.6428
$ws.Range("K58").Value = 574.931
$ws.Range("L58").Value = 2319.6428
$ws.Range("M58").Value = -371.931
$ws.Range("N58").Value = -2725.6428
$ws.Range("H99").Value = 15629500
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("H122").Value = 1160.1
$ws.Range("J122").Value = 1377.8
$ws.Range("L122").Value = 4133.4
$ws.Range("N122").Value = -9033.4
$ws.Range("H126").Value = 15629500
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 3308
$ws.Range("I132").Value = 2733.3333
$ws.Range("J132").Value = 3997.6
$ws.Range("K132").Value = 8199.999899999999
$ws.Range("L132").Value = 11992.8
$ws.Range("M132").Value = -5669.999899999999
$ws.Range("N132").Value = -17052.8
$ws.Range("H136").Value = 1142.9767
$ws.Range("I136").Value = 574.931
$ws.Range("J136").Value = 2319.6428
$ws.Range("K136").Value = 1724.793
$ws.Range("L136").Value = 6958.928400000001
$ws.Range("M136").Value = 825.2069999999999
$ws.Range("N136").Value = -12058.9284
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 855.4
$ws.Range("I23").Value = 1359.6666
$ws.Range("J23").Value = 442.81818
$ws.Range("K23").Value = 4078.9998
$ws.Range("L23").Value = 1328.45454
$ws.Range("M23").Value = -3843.9998
$ws.Range("N23").Value = -1798.45454
$ws.Range("H131").Value = 2139.2449
$ws.Range("I131").Value = 366.66666
$ws.Range("J131").Value = 2386.5813
$ws.Range("K131").Value = 1099.99998
$ws.Range("L131").Value = 7159.743899999999
$ws.Range("M131").Value = 3940.00002
$ws.Range("N131").Value = -17239.7439
$ws.Range("H141").Value = 4513.75
$ws.Range("J141").Value = 2250
$ws.Range("L141").Value = 6750
$ws.Range("N141").Value = -17110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2156.1177
$ws.Range("I102").Value = 2039.25
$ws.Range("J102").Value = 2260
$ws.Range("K102").Value = 2039.25
$ws.Range("L102").Value = 2260
$ws.Range("M102").Value = -417.25
$ws.Range("N102").Value = -5504
$ws.Range("H121").Value = 29000
$ws.Range("J121").Value = 29000
$ws.Range("L121").Value = 29000
$ws.Range("N121").Value = -32494
$ws.Range("H132").Value = 2429.5112
$ws.Range("I132").Value = 1934.1428
$ws.Range("J132").Value = 4163.3
$ws.Range("K132").Value = 5802.428400000001
$ws.Range("L132").Value = 12489.9
$ws.Range("M132").Value = -3272.428400000001
$ws.Range("N132").Value = -17549.9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1346.4
$ws.Range("I16").Value = 1357.5
$ws.Range("J16").Value = 1302
$ws.Range("K16").Value = 1357.5
$ws.Range("L16").Value = 1302
$ws.Range("M16").Value = -1187.5
$ws.Range("N16").Value = -1642
$ws.Range("H40").Value = 3055.5
$ws.Range("I40").Value = 2244.2856
$ws.Range("J40").Value = 3492.3076
$ws.Range("K40").Value = 2244.2856
$ws.Range("L40").Value = 3492.3076
$ws.Range("M40").Value = -2108.2856
$ws.Range("N40").Value = -3764.3076
$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 45000
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45450
$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 45000
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46560
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 26029.334
$ws.Range("J121").Value = 26029.334
$ws.Range("L121").Value = 26029.334
$ws.Range("N121").Value = -29523.334
$ws.Range("H122").Value = 101367.8
$ws.Range("I122").Value = 126336
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 379008
$ws.Range("L122").Value = 4485
$ws.Range("M122").Value = -376558
$ws.Range("N122").Value = -9385
$ws.Range("H124").Value = 39933.332
$ws.Range("J124").Value = 39933.332
$ws.Range("L124").Value = 39933.332
$ws.Range("N124").Value = -49753.332
$ws.Range("H126").Value = 62861.062
$ws.Range("I126").Value = 77282.84
$ws.Range("J126").Value = 366.66666
$ws.Range("K126").Value = 231848.52
$ws.Range("L126").Value = 1099.99998
$ws.Range("M126").Value = -229378.52
$ws.Range("N126").Value = -6039.999980000001
$ws.Range("H135").Value = 86611
$ws.Range("J135").Value = 86611
$ws.Range("L135").Value = 86611
$ws.Range("N135").Value = -96751
